$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values per the daily price refresh (Day advanced by one, new hourly prices)
$ws.Range("A2").Value = 45991
$ws.Range("B2").Value = 76.84999999999999
$ws.Range("C2").Value = 70.76000000000001
$ws.Range("D2").Value = 68.5
$ws.Range("E2").Value = 66.40000000000001
$ws.Range("F2").Value = 65.95999999999999
$ws.Range("G2").Value = 66.45
$ws.Range("H2").Value = 65.20999999999999
$ws.Range("I2").Value = 69.37
$ws.Range("J2").Value = 73.84
$ws.Range("K2").Value = 77.83
$ws.Range("L2").Value = 67.92
$ws.Range("M2").Value = 56.16
$ws.Range("N2").Value = 41.55
$ws.Range("O2").Value = 48.8
$ws.Range("P2").Value = 51.49
$ws.Range("Q2").Value = 60.44
$ws.Range("R2").Value = 80.34999999999999
$ws.Range("S2").Value = 99.27
$ws.Range("T2").Value = 105.38
$ws.Range("U2").Value = 106.81
$ws.Range("V2").Value = 101.42
$ws.Range("W2").Value = 95.93000000000001
$ws.Range("X2").Value = 93.59999999999999
$ws.Range("Y2").Value = 87.70999999999999
$ws.Range("Z2").Value = 74.92
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 97.95
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 106.1
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 98.68000000000001
$ws.Range("AG2").Value = "1h-15h"
